# Actualización automática 2025-12-04 16:30:07
#
# Insert a new salesperson row ("PARRALES LIRIANO CARLOS JULIO") right
# before "PERDOMO BRIONES JOSÉ ALBERTO" (alphabetical order) in both the
# "VENTAS POR GRUPO" and "VENTA MENSUAL" sheets. The new row becomes row
# 45 on each sheet, pushing every following row (including the trailing
# totals / "X de Y" summary row) down by one.

$wb = $excel.ActiveWorkbook

$sheetNames = @("VENTAS POR GRUPO", "VENTA MENSUAL")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Determine the used range extent before the insert.
    $usedRange = $ws.UsedRange
    $lastCol = $usedRange.Columns.Count
    $lastRowBefore = $usedRange.Rows.Count

    # Push row 45 (and everything below it) down by one row.
    $ws.Rows("45:45").Insert()

    # Fill in the new row's data.
    $ws.Cells.Item(45, 1).Value2 = "OFICINA-CATAECSA"
    $ws.Cells.Item(45, 2).Value2 = "PARRALES LIRIANO CARLOS JULIO"

    for ($col = 3; $col -le $lastCol; $col++) {
        $ws.Cells.Item(45, $col).Value2 = 0
    }

    # The trailing summary row (e.g. "0 de 56") counts out of the total
    # number of advisors on the sheet, expressed as literal text. It
    # doesn't recompute on its own, so bump "de 56" -> "de 57" now that
    # one more advisor row exists. It now lives one row further down
    # because of the insert above.
    $summaryRow = $lastRowBefore + 1
    for ($col = 3; $col -le $lastCol; $col++) {
        $cell = $ws.Cells.Item($summaryRow, $col)
        $current = $cell.Value2
        if ($current -ne $null -and $current -like "*de 56*") {
            $cell.Value2 = $current -replace "de 56", "de 57"
        }
    }
}
